$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '66.250.92'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -1.81%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.272.16'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -2.10%  '

$ws.Range("E4").Value = '  +0.00%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '578.00'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '178.88'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -2.40%  '

$ws.Range("E7").Value = '  +3.72%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -2.59%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '6.72'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +0.77%  '

$ws.Range("E11").Value = '  -1.52%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '3.849.09'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -1.88%  '

$ws.Range("E13").Value = '  -3.78%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '66.288.80'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -2.03%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '26.33'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -3.81%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.296.46'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -1.05%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.0000163'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -2.44%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '434.37'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.25%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '5.50'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -2.95%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '13.15'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.10%  '

$ws.Range("E21").Value = '  -4.55%  '

$ws.Range("E22").Value = '  -2.91%  '

$ws.Range("E23").Value = '  +0.14%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '3.421.44'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.60%  '

$ws.Range("E25").Value = '  -1.64%  '

$ws.Range("E26").Value = '  +2.51%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.0000112'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -5.98%  '

$ws.Range("E28").Value = '  -2.30%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("E30").Value = '  -1.91%  '

$ws.Range("E31").Value = '  -3.04%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.16'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.30%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '6.56'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.29%  '

$ws.Range("E35").Value = '  -4.05%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '157.15'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.93%  '

$ws.Range("E37").Value = '  -4.91%  '

$ws.Range("E38").Value = '  -2.62%  '

$ws.Range("E39").Value = '  -3.42%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.764.52'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -2.08%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.771'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.48%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '4.28'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -4.10%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '40.27'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '6.02'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -3.28%  '

$ws.Range("E45").Value = '  -2.34%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '318.82'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -1.85%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '2.28'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -3.74%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '23.11'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -6.22%  '

$ws.Range("E49").Value = '  -2.54%  '

$ws.Range("E50").Value = '  +2.30%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.06%  '
